# Start mit Auswertung nach Einrichtungen
# Inserts two new "Stiftung Jugendhilfe aktiv Region ..." rows into the
# Einrichtungen overview and refreshes the SBBZ/HZE counts for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final state of the data rows (Einrichtung, Leitung, SBBZ, HZE), in row order.
$rows = @(
    @('Jugendhilfeverbund Kinderheim Rodt', 3, 9, 7),
    @('Ev Jugendhilfe Friedenshort', 6, 4, 21),
    @('Jugendhilfe Hoffmannhaus Wilhelmsdorf', 6, 6, 10),
    @('Stiftung Jugendhilfe aktiv Region Stuttgart', 4, 4, 6),
    @('Jugendhilfe Hochdorf', 4, 0, 10),
    @('Sophienpflege Tübingen', 2, 2, 6),
    @('Stiftung Tragwerk KirchheimTeck', 7, 2, 23),
    @('Stiftung Jugendhilfe aktiv Region Böblingen', 4, 2, 11),
    @('Mutpol Tuttlingen', 6, 24, 28),
    @('Stiftung Jugendhilfe aktiv Region Esslingen', 4, 0, 3),
    @('eva Heidenheim', 2, 0, 0),
    @('Jugendhilfe Deggingen', 2, 2, 2),
    @('Scout am Löwentor Stuttgart', 2, 0, 7),
    @('Jugendhilfe Korntal', 5, 1, 7),
    @('Johannes_Falk_Haus Stuttgart', 0, 0, 2),
    @('Oberlin Jugendhilfe Reutlingen', 2, 1, 6),
    @('Jugendhilfeverbund Paulinenpflege Winnenden', 6, 0, 22),
    @('Diakonische Jugendhilfe Heilbronn (DJHN)', 3, 2, 11),
    @('eva Stuttgart', 3, 0, 42),
    @('Mariaberger Heime Gammertingen', 2, 1, 4),
    @('Oberlin Ulm', 8, 0, 8),
    @('Diasporahaus Bietenhausen', 5, 15, 10),
    @('Weraheim Hebsack', 1, 0, 0),
    @('Kinder und Jugendhilfe Karlshöhe Ludwigsburg', 2, 0, 4)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}
